# Write-Log design goals.pptx - add "Comment based help" as a new bullet
# point at the end of the Content Placeholder on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Content Placeholder 2" is the second shape on the slide (the first is
# the title). Grab its text range and append a new paragraph at the end.
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Use InsertAfter (rather than re-assigning .Text with the existing
# contents) so the untouched runs keep their exact original characters
# (e.g. curly quotes) instead of being round-tripped/normalised.
# "`r" starts a brand-new paragraph (like pressing Enter in PowerPoint),
# as opposed to a soft line-break within the same paragraph.
$tr.InsertAfter("`rComment based help") | Out-Null
